$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this pushes the existing rows 48-89
# down to 49-90 (and the sheet dimension grows from T89 to T90).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with a new weekly Mango price record
# (same market/product metadata as its neighbours, new date & values).
$ws.Range("A48").Value = 5
$ws.Range("B48").Value = 'Macroferia Regional de Talca'
$ws.Range("C48").Value = 'Maule'
$ws.Range("D48").Value = 44512
$ws.Range("E48").Value = 7
$ws.Range("F48").Value = 'Fruta'
$ws.Range("G48").Value = 100108
$ws.Range("H48").Value = 'Tropicales y subtropicales'
$ws.Range("I48").Value = 100108002
$ws.Range("J48").Value = 'Mango'
$ws.Range("K48").Value = 'Sin especificar'
$ws.Range("L48").Value = 'Primera'
$ws.Range("M48").Value = 200
$ws.Range("N48").Value = 6000
$ws.Range("O48").Value = 6000
$ws.Range("P48").Value = 6000
$ws.Range("Q48").Value = '$/bandeja 4 kilos'
$ws.Range("R48").Value = 'Perú'
$ws.Range("S48").Value = 1500
$ws.Range("T48").Value = 4
